$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A470").Value = "Edit"
$ws.Range("B470").Value = 26662
$ws.Range("C470").Value = "Artiklar och tjänster"
$ws.Range("D470").Value = "Transaktionstyp"

$ws.Range("A471").Value = "SafGrid"
$ws.Range("B471").Value = 21530
$ws.Range("C471").Value = "Artiklar och tjänster"
$ws.Range("D471").Value = "Stafflingar"

$ws.Range("A472").Value = "SafGrid"
$ws.Range("B472").Value = 21521
$ws.Range("C472").Value = "Artiklar och tjänster"
$ws.Range("D472").Value = "Priser"

$ws.Range("A473").Value = "Edit"
$ws.Range("B473").Value = 22206
$ws.Range("C473").Value = "Hämta projekt"
$ws.Range("D473").Value = "Projektnummer"

$ws.Range("A474").Value = "Edit"
$ws.Range("B474").Value = 22207
$ws.Range("C474").Value = "Hämta projekt"
$ws.Range("D474").Value = "Namn"

$ws.Range("A475").Value = "Edit"
$ws.Range("B475").Value = 22208
$ws.Range("C475").Value = "Hämta projekt"
$ws.Range("D475").Value = "Projektet startar"

$ws.Range("A476").Value = "Edit"
$ws.Range("B476").Value = 22209
$ws.Range("C476").Value = "Hämta projekt"
$ws.Range("D476").Value = "Projektet slutar"

$ws.Range("A477").Value = "Edit"
$ws.Range("B477").Value = 22212
$ws.Range("C477").Value = "Hämta projekt"
$ws.Range("D477").Value = "Kundnummer"

$ws.Range("A478").Value = "Edit"
$ws.Range("B478").Value = 22213
$ws.Range("C478").Value = "Hämta projekt"
$ws.Range("D478").Value = "Kundnamn"

$ws.Range("A479").Value = "Edit"
$ws.Range("B479").Value = 22215
$ws.Range("C479").Value = "Hämta projekt"
$ws.Range("D479").Value = "Kundens ordernummer"

$ws.Range("A480").Value = "Edit"
$ws.Range("B480").Value = 22216
$ws.Range("C480").Value = "Hämta projekt"
$ws.Range("D480").Value = "Kontaktperson"

$ws.Range("A481").Value = "Edit"
$ws.Range("B481").Value = 22217
$ws.Range("C481").Value = "Hämta projekt"
$ws.Range("D481").Value = "Telefon"

$ws.Range("A482").Value = "Edit"
$ws.Range("B482").Value = 22218
$ws.Range("C482").Value = "Hämta projekt"
$ws.Range("D482").Value = "Fax"

$ws.Range("A483").Value = "Edit"
$ws.Range("B483").Value = 22219
$ws.Range("C483").Value = "Hämta projekt"
$ws.Range("D483").Value = "Arbetsplats"

$ws.Range("A484").Value = "Edit"
$ws.Range("B484").Value = 22220
$ws.Range("C484").Value = "Hämta projekt"
$ws.Range("D484").Value = "Arbetsplats, fortsättning"

$ws.Range("A485").Value = "Edit"
$ws.Range("B485").Value = 22214
$ws.Range("C485").Value = "Hämta projekt"
$ws.Range("D485").Value = "Vårt ordernummer"

$ws.Range("A486").Value = "ComboBox"
$ws.Range("B486").Value = 22225
$ws.Range("C486").Value = "Hämta projekt"
$ws.Range("D486").Value = "Kontoplanstyp"

$ws.Range("A487").Value = "Edit"
$ws.Range("B487").Value = 22211
$ws.Range("C487").Value = "Hämta projekt"
$ws.Range("D487").Value = "Ansvarig"

$ws.Range("A488").Value = "Edit"
$ws.Range("B488").Value = 22221
$ws.Range("C488").Value = "Hämta projekt"
$ws.Range("D488").Value = "Anteckning 1"

$ws.Range("A489").Value = "Edit"
$ws.Range("B489").Value = 22222
$ws.Range("C489").Value = "Hämta projekt"
$ws.Range("D489").Value = "Anteckning 2"

$ws.Range("A490").Value = "Edit"
$ws.Range("B490").Value = 22223
$ws.Range("C490").Value = "Hämta projekt"
$ws.Range("D490").Value = "Anteckning 3"

$ws.Range("A491").Value = "Edit"
$ws.Range("B491").Value = 22224
$ws.Range("C491").Value = "Hämta projekt"
$ws.Range("D491").Value = "Anteckning 4"

$ws.Range("B482").Select()
